$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.757.95"
$ws.Range("E2").Value = '  -2.66%  '

$ws.Range("D3").Value = "'1.742.89"
$ws.Range("E3").Value = '  -5.15%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = "'237.44"
$ws.Range("E5").Value = '  -9.23%  '

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = "'0.5029"
$ws.Range("E7").Value = '  -6.64%  '

$ws.Range("D8").Value = "'41.94"
$ws.Range("E8").Value = '  -6.58%  '

$ws.Range("D9").Value = "'0.2673"
$ws.Range("E9").Value = '  -11.34%  '

$ws.Range("D10").Value = "'0.06137"
$ws.Range("E10").Value = '  -10.83%  '

$ws.Range("D11").Value = "'1.745.19"
$ws.Range("E11").Value = '  -4.99%  '

$ws.Range("D12").Value = "'0.06922"
$ws.Range("E12").Value = '  -3.30%  '

$ws.Range("D13").Value = "'15.39"
$ws.Range("E13").Value = '  -12.95%  '

$ws.Range("D14").Value = "'4.511"
$ws.Range("E14").Value = '  -9.63%  '

$ws.Range("D15").Value = "'0.5977"
$ws.Range("E15").Value = '  -19.01%  '

$ws.Range("D16").Value = "'76.73"
$ws.Range("E16").Value = '  -13.97%  '

$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").Value = "'25.766.36"
$ws.Range("E19").Value = '  -2.73%  '

$ws.Range("E20").Value = '  -13.49%  '

$ws.Range("D21").Value = "'11.55"
$ws.Range("E21").Value = '  -16.55%  '

$ws.Range("D22").Value = "'1.966.27"
$ws.Range("E22").Value = '  -5.33%  '

$ws.Range("D23").Value = "'4.046"
$ws.Range("E23").Value = '  -11.86%  '

$ws.Range("D24").Value = "'5.219"
$ws.Range("E24").Value = '  -12.72%  '

$ws.Range("D25").Value = "'8.135"

$ws.Range("D26").Value = "'137.32"
$ws.Range("E26").Value = '  -3.92%  '

$ws.Range("D27").Value = "'1.516"
$ws.Range("E27").Value = '  -10.54%  '

$ws.Range("D28").Value = "'14.97"
$ws.Range("E28").Value = '  -11.79%  '

$ws.Range("D29").Value = "'1.805"
$ws.Range("E29").Value = '  -17.77%  '

$ws.Range("D30").Value = "'104.04"
$ws.Range("E30").Value = '  -5.99%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'0.08109"
$ws.Range("E31").Value = '  -8.11%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = "'3.748"
$ws.Range("E32").Value = '  -11.55%  '

$ws.Range("D33").Value = "'3.459"
$ws.Range("E33").Value = '  -14.24%  '

$ws.Range("D34").Value = "'0.04525"
$ws.Range("E34").Value = '  -6.08%  '

$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").Value = "'2.633"
$ws.Range("E36").Value = '  -9.77%  '

$ws.Range("D37").Value = "'0.9773"
$ws.Range("E37").Value = '  -13.57%  '

$ws.Range("D38").Value = "'0.6095"
$ws.Range("E38").Value = '  -16.49%  '

$ws.Range("D39").Value = "'2.667"
$ws.Range("E39").Value = '  -13.77%  '

$ws.Range("D40").Value = "'0.01552"
$ws.Range("E40").Value = '  -9.48%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'1.906"
$ws.Range("E42").Value = '  -15.71%  '

$ws.Range("D43").Value = "'101.74"
$ws.Range("E43").Value = '  -5.74%  '

$ws.Range("D44").Value = "'0.3803"
$ws.Range("E44").Value = '  -19.43%  '

$ws.Range("D45").Value = "'5.068"
$ws.Range("E45").Value = '  -13.96%  '

$ws.Range("D46").Value = "'0.7327"
$ws.Range("E46").Value = '  -19.07%  '

$ws.Range("D47").Value = "'0.05370"
$ws.Range("E47").Value = '  -6.98%  '

$ws.Range("D48").Value = "'0.1105"
$ws.Range("E48").Value = '  -11.19%  '

$ws.Range("D49").Value = "'30.14"
$ws.Range("E49").Value = '  -13.45%  '

$ws.Range("D50").Value = "'5.905"
$ws.Range("E50").Value = '  -20.05%  '

$ws.Range("D51").Value = "'52.51"
$ws.Range("E51").Value = '  -12.52%  '
